$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update stats for 2025-12 (row 25)
$ws.Range("B25").Value = 6481
$ws.Range("D25").Value = 6044502
$ws.Range("E25").Value = 932.6495911124827
$ws.Range("F25").Value = 10.01527754201323
$ws.Range("H25").Value = 26.59051202907905
